# Weekly update: a new price-record row is inserted at row 37 ("Hortaliza,
# Femacal de La Calera - Haba"). All existing data rows from 37 downward
# shift down by one row; the freshly inserted row 37 is then populated
# with the new week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 37, pushing rows 37:54
# down to 38:55.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new record.
$ws.Cells.Item(37, 1).Value  = 3
$ws.Cells.Item(37, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(37, 3).Value  = "Coquimbo"
$ws.Cells.Item(37, 4).Value  = 44452
$ws.Cells.Item(37, 5).Value  = 5
$ws.Cells.Item(37, 6).Value  = 100112026
$ws.Cells.Item(37, 7).Value  = "Haba"
$ws.Cells.Item(37, 8).Value  = "Sin especificar"
$ws.Cells.Item(37, 9).Value  = "Primera"
$ws.Cells.Item(37, 10).Value = 35
$ws.Cells.Item(37, 11).Value = 15000
$ws.Cells.Item(37, 12).Value = 15000
$ws.Cells.Item(37, 13).Value = 15000
$ws.Cells.Item(37, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(37, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(37, 16).Value = 600
$ws.Cells.Item(37, 17).Value = 25
$ws.Cells.Item(37, 18).Value = "Hortaliza"
